# "Pais" sheet refresh (new scrape snapshot):
#   - A1 "datos actualizados" timestamp bumped to 15:06
#   - Updated Casos totales/Nuevos casos/Casos activos/Recuperados/
#     Casos criticos/Muertes hoy/Muertes for the affected countries
#   - Because the sheet is ranked by "Casos totales" desc, a few
#     neighbouring countries swap rows as their totals cross over
#     (Colombia/Paises Bajos, Gabon/Kenia, Etiopia/Bulgaria,
#     Guayana Francesa/Sierra Leona, Groenlandia/Islas Malvinas,
#     Montserrat/Seychelles, Papua Nueva Guinea/Islas Virgenes Britanicas)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 14 de Junio de 2020 a las 15:06'

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 2142540
$ws.Cells.Item(4, 3).Value = 316
$ws.Cells.Item(4, 5).Value = 1170901

# Row 7: India
$ws.Cells.Item(7, 2).Value = 322777
$ws.Cells.Item(7, 3).Value = 1151
$ws.Cells.Item(7, 4).Value = 163019
$ws.Cells.Item(7, 5).Value = 150552
$ws.Cells.Item(7, 7).Value = 7
$ws.Cells.Item(7, 8).Value = 9206

# Row 19: Arabia Saudita
$ws.Cells.Item(19, 2).Value = 127541
$ws.Cells.Item(19, 3).Value = 4233
$ws.Cells.Item(19, 4).Value = 84720
$ws.Cells.Item(19, 5).Value = 41849
$ws.Cells.Item(19, 7).Value = 40
$ws.Cells.Item(19, 8).Value = 972

# Row 23: Catar
$ws.Cells.Item(23, 2).Value = 79602
$ws.Cells.Item(23, 3).Value = 1186
$ws.Cells.Item(23, 4).Value = 56898
$ws.Cells.Item(23, 5).Value = 22631
$ws.Cells.Item(23, 7).Value = 3
$ws.Cells.Item(23, 8).Value = 73

# Row 27: Suecia
$ws.Cells.Item(27, 2).Value = 51614
$ws.Cells.Item(27, 3).Value = 38

# Row 28: Paises Bajos
$ws.Cells.Item(28, 1).Value = 'Paises Bajos'
$ws.Cells.Item(28, 2).Value = 48783
$ws.Cells.Item(28, 3).Value = 143
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 7).Value = 2
$ws.Cells.Item(28, 8).Value = 6059

# Row 29: Colombia
$ws.Cells.Item(29, 1).Value = 'Colombia'
$ws.Cells.Item(29, 2).Value = 48746
$ws.Cells.Item(29, 4).Value = 19426
$ws.Cells.Item(29, 5).Value = 27728
$ws.Cells.Item(29, 8).Value = 1592

# Row 35: Portugal
$ws.Cells.Item(35, 2).Value = 36690
$ws.Cells.Item(35, 3).Value = 227
$ws.Cells.Item(35, 4).Value = 22669
$ws.Cells.Item(35, 5).Value = 12504
$ws.Cells.Item(35, 7).Value = 5
$ws.Cells.Item(35, 8).Value = 1517

# Row 39: Argentina
$ws.Cells.Item(39, 4).Value = 9564
$ws.Cells.Item(39, 5).Value = 19912
$ws.Cells.Item(39, 7).Value = 4
$ws.Cells.Item(39, 8).Value = 819

# Row 53: Austria
$ws.Cells.Item(53, 2).Value = 17109
$ws.Cells.Item(53, 3).Value = 31
$ws.Cells.Item(53, 4).Value = 16059
$ws.Cells.Item(53, 5).Value = 373

# Row 57: Serbia
$ws.Cells.Item(57, 2).Value = 12310
$ws.Cells.Item(57, 3).Value = 59
$ws.Cells.Item(57, 5).Value = 708
$ws.Cells.Item(57, 7).Value = 1
$ws.Cells.Item(57, 8).Value = 254

# Row 58: Dinamarca
$ws.Cells.Item(58, 2).Value = 12193
$ws.Cells.Item(58, 3).Value = 54
$ws.Cells.Item(58, 4).Value = 11068
$ws.Cells.Item(58, 5).Value = 528

# Row 76: Uzbekistan
$ws.Cells.Item(76, 2).Value = 5051
$ws.Cells.Item(76, 3).Value = 85
$ws.Cells.Item(76, 5).Value = 1122

# Row 87: Kenia
$ws.Cells.Item(87, 1).Value = 'Kenia'
$ws.Cells.Item(87, 2).Value = 3594
$ws.Cells.Item(87, 3).Value = 137
$ws.Cells.Item(87, 4).Value = 1221
$ws.Cells.Item(87, 5).Value = 2273
$ws.Cells.Item(87, 8).Value = 100

# Row 88: Gabon
$ws.Cells.Item(88, 1).Value = 'Gabon'
$ws.Cells.Item(88, 2).Value = 3463
$ws.Cells.Item(88, 4).Value = 1024
$ws.Cells.Item(88, 5).Value = 2416
$ws.Cells.Item(88, 8).Value = 23

# Row 89: Etiopia
$ws.Cells.Item(89, 1).Value = 'Etiopia'
$ws.Cells.Item(89, 2).Value = 3345
$ws.Cells.Item(89, 3).Value = 179
$ws.Cells.Item(89, 4).Value = 545
$ws.Cells.Item(89, 5).Value = 2743
$ws.Cells.Item(89, 7).Value = 2
$ws.Cells.Item(89, 8).Value = 57

# Row 90: Bulgaria
$ws.Cells.Item(90, 1).Value = 'Bulgaria'
$ws.Cells.Item(90, 2).Value = 3266
$ws.Cells.Item(90, 4).Value = 1723
$ws.Cells.Item(90, 5).Value = 1371
$ws.Cells.Item(90, 8).Value = 172

# Row 104: Islandia
$ws.Cells.Item(104, 2).Value = 1810
$ws.Cells.Item(104, 3).Value = 2
$ws.Cells.Item(104, 4).Value = 1796

# Row 121: Sierra Leona
$ws.Cells.Item(121, 1).Value = 'Sierra Leona'
$ws.Cells.Item(121, 2).Value = 1169
$ws.Cells.Item(121, 3).Value = 37
$ws.Cells.Item(121, 4).Value = 680
$ws.Cells.Item(121, 5).Value = 438
$ws.Cells.Item(121, 8).Value = 51

# Row 122: Guayana Francesa
$ws.Cells.Item(122, 1).Value = 'Guayana Francesa'
$ws.Cells.Item(122, 2).Value = 1161
$ws.Cells.Item(122, 4).Value = 520
$ws.Cells.Item(122, 5).Value = 639
$ws.Cells.Item(122, 8).Value = 2

# Row 125: Tunez
$ws.Cells.Item(125, 2).Value = 1096
$ws.Cells.Item(125, 3).Value = 2
$ws.Cells.Item(125, 4).Value = 998
$ws.Cells.Item(125, 5).Value = 49

# Row 129: Burkina Faso
$ws.Cells.Item(129, 2).Value = 894
$ws.Cells.Item(129, 3).Value = 2
$ws.Cells.Item(129, 5).Value = 42

# Row 143: Mozambique
$ws.Cells.Item(143, 2).Value = 583
$ws.Cells.Item(143, 3).Value = 30
$ws.Cells.Item(143, 5).Value = 429
$ws.Cells.Item(143, 7).Value = 1
$ws.Cells.Item(143, 8).Value = 3

# Row 206: Islas Malvinas
$ws.Cells.Item(206, 1).Value = 'Islas Malvinas'

# Row 207: Groenlandia
$ws.Cells.Item(207, 1).Value = 'Groenlandia'

# Row 210: Seychelles
$ws.Cells.Item(210, 1).Value = 'Seychelles'
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 8).Value = 0

# Row 211: Montserrat
$ws.Cells.Item(211, 1).Value = 'Montserrat'
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 8).Value = 1

# Row 213: Islas Virgenes Britanicas
$ws.Cells.Item(213, 1).Value = 'Islas Virgenes Britanicas'
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1

# Row 214: Papua Nueva Guinea
$ws.Cells.Item(214, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0
